$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Append the new log row (row 14)
$ws.Range("A14").Value = "Hebben jullie toevallig al iets gehoord?"
$ws.Range("B14").Value = "mailmind.test@zohomail.eu"
$ws.Range("C14").Value = "Testmail #1: Hebben jullie toevallig al iets gehoord?"
$ws.Range("D14").Value = "Klantenservice / Contact"
$ws.Range("E14").Value = "Bedankt, we hebben dit doorgestuurd naar klantenservice@bedrijf.nl."
$ws.Range("F14").Value = "2025-08-06 20:06:49"
$ws.Range("G14").Value = "Ja"
$ws.Range("H14").Value = "Ja"
$ws.Range("I14").Value = "Nee"
$ws.Range("J14").Value = "Nee"

# Extend the conditional formatting ranges to include the new row
$colsOld = @("D2:D13", "G2:G13", "H2:H13", "I2:I13", "J2:J13")
$colsNew = @("D2:D14", "G2:G14", "H2:H14", "I2:I14", "J2:J14")
for ($i = 0; $i -lt $colsOld.Count; $i++) {
    $fcs = $ws.Range($colsOld[$i]).FormatConditions
    $newRange = $ws.Range($colsNew[$i])
    for ($j = 1; $j -le $fcs.Count(); $j++) {
        $fcs.Item($j).ModifyAppliesToRange($newRange)
    }
}

# Update the Dashboard summary count for "Klantenservice / Contact"
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B4").Value = 2
